$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 168.96428
$ws.Range("I33").Value = 176.44
$ws.Range("J33").Value = 106.666664
$ws.Range("K33").Value = 176.44
$ws.Range("L33").Value = 106.666664
$ws.Range("M33").Value = 52.56
$ws.Range("N33").Value = -564.666664

$ws.Range("H74").Value = 8932949
$ws.Range("I74").Value = 4184.2856
$ws.Range("J74").Value = 17861714
$ws.Range("K74").Value = 4184.2856
$ws.Range("L74").Value = 17861714
$ws.Range("M74").Value = -3248.2856
$ws.Range("N74").Value = -17863586

$ws.Range("H77").Value = 8932949
$ws.Range("I77").Value = 4184.2856
$ws.Range("J77").Value = 17861714
$ws.Range("K77").Value = 20921.428
$ws.Range("L77").Value = 89308570
$ws.Range("M77").Value = -16241.428
$ws.Range("N77").Value = -89317930

$ws.Range("H116").Value = 27783612
$ws.Range("I116").Value = 125001750
$ws.Range("K116").Value = 125001750
$ws.Range("M116").Value = -124998308

$ws.Range("H137").Value = 123154.336
$ws.Range("I137").Value = 149451.48
$ws.Range("K137").Value = 448354.4400000001
$ws.Range("M137").Value = -445804.4400000001

$ws.Range("H138").Value = 2821.7175
$ws.Range("I138").Value = 2343.3635
$ws.Range("J138").Value = 2892.8242
$ws.Range("K138").Value = 7030.0905
$ws.Range("L138").Value = 8678.472600000001
$ws.Range("M138").Value = -1890.0905
$ws.Range("N138").Value = -18958.4726

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1084.375
$ws.Range("I2").Value = 1215
$ws.Range("J2").Value = 866.6667
$ws.Range("K2").Value = 1215
$ws.Range("L2").Value = 866.6667
$ws.Range("M2").Value = -1102
$ws.Range("N2").Value = -1092.6667

$ws.Range("H32").Value = 8081.0234
$ws.Range("I32").Value = 5906.586
$ws.Range("J32").Value = 17594.188
$ws.Range("K32").Value = 5906.586
$ws.Range("L32").Value = 17594.188
$ws.Range("M32").Value = -5619.586
$ws.Range("N32").Value = -18168.188

$ws.Range("H45").Value = 2438.0967
$ws.Range("I45").Value = 2258.6316
$ws.Range("K45").Value = 2258.6316
$ws.Range("M45").Value = -1881.6316

$ws.Range("H74").Value = 29413270
$ws.Range("I74").Value = 38462160
$ws.Range("K74").Value = 38462160
$ws.Range("M74").Value = -38461286

$ws.Range("H77").Value = 29413270
$ws.Range("I77").Value = 38462160
$ws.Range("K77").Value = 192310800
$ws.Range("M77").Value = -192306432

$ws.Range("H110").Value = 974.75
$ws.Range("I110").Value = 906.7059
$ws.Range("J110").Value = 1140
$ws.Range("K110").Value = 906.7059
$ws.Range("L110").Value = 1140
$ws.Range("M110").Value = 1138.2941
$ws.Range("N110").Value = -5230

$ws.Range("H116").Value = 1084.375
$ws.Range("I116").Value = 1215
$ws.Range("J116").Value = 866.6667
$ws.Range("K116").Value = 1215
$ws.Range("L116").Value = 866.6667
$ws.Range("M116").Value = 1079
$ws.Range("N116").Value = -5454.6667

$ws.Range("H132").Value = 13235.467
$ws.Range("I132").Value = 1726.5555
$ws.Range("J132").Value = 59271.11
$ws.Range("K132").Value = 5179.666499999999
$ws.Range("L132").Value = 177813.33
$ws.Range("M132").Value = -2649.666499999999
$ws.Range("N132").Value = -182873.33

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1084.375
$ws.Range("I3").Value = 1215
$ws.Range("J3").Value = 866.6667
$ws.Range("K3").Value = 1215
$ws.Range("L3").Value = 866.6667
$ws.Range("M3").Value = -1101
$ws.Range("N3").Value = -1094.6667

$ws.Range("H32").Value = 6013
$ws.Range("I32").Value = 2026
$ws.Range("J32").Value = 10000
$ws.Range("K32").Value = 2026
$ws.Range("L32").Value = 10000
$ws.Range("M32").Value = -1642
$ws.Range("N32").Value = -10768

$ws.Range("H86").Value = 1947.0588
$ws.Range("I86").Value = 1725
$ws.Range("K86").Value = 1725
$ws.Range("M86").Value = -602

$ws.Range("H89").Value = 1947.0588
$ws.Range("I89").Value = 1725
$ws.Range("K89").Value = 8625
$ws.Range("M89").Value = -3009

$ws.Range("H107").Value = 1379.0303
$ws.Range("I107").Value = 1083.6
$ws.Range("K107").Value = 1083.6
$ws.Range("M107").Value = 836.4000000000001

$ws.Range("H134").Value = 3409.9744
$ws.Range("I134").Value = 3447.0789
$ws.Range("K134").Value = 10341.2367
$ws.Range("M134").Value = -7806.236699999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4608.171
$ws.Range("I31").Value = 2352.389
$ws.Range("J31").Value = 6373.5654
$ws.Range("K31").Value = 2352.389
$ws.Range("L31").Value = 6373.5654
$ws.Range("M31").Value = -2057.389
$ws.Range("N31").Value = -6963.5654

$ws.Range("H34").Value = 4608.171
$ws.Range("I34").Value = 2352.389
$ws.Range("J34").Value = 6373.5654
$ws.Range("K34").Value = 2352.389
$ws.Range("L34").Value = 6373.5654
$ws.Range("M34").Value = -2150.389
$ws.Range("N34").Value = -6777.5654

$ws.Range("H58").Value = 40591.848
$ws.Range("I58").Value = 2414.5715
$ws.Range("J58").Value = 85132
$ws.Range("K58").Value = 2414.5715
$ws.Range("L58").Value = 85132
$ws.Range("M58").Value = -2211.5715
$ws.Range("N58").Value = -85538

$ws.Range("H86").Value = 16683723
$ws.Range("I86").Value = 2500
$ws.Range("J86").Value = 18537192
$ws.Range("K86").Value = 2500
$ws.Range("L86").Value = 18537192
$ws.Range("M86").Value = -1377
$ws.Range("N86").Value = -18539438

$ws.Range("H89").Value = 16683723
$ws.Range("I89").Value = 2500
$ws.Range("J89").Value = 18537192
$ws.Range("K89").Value = 12500
$ws.Range("L89").Value = 92685960
$ws.Range("M89").Value = -6884
$ws.Range("N89").Value = -92697192

$ws.Range("H131").Value = 39325.16
$ws.Range("J131").Value = 39325.16
$ws.Range("L131").Value = 39325.16
$ws.Range("N131").Value = -49405.16

$ws.Range("H132").Value = 3320.6667
$ws.Range("I132").Value = 2213.875
$ws.Range("J132").Value = 6862.4
$ws.Range("K132").Value = 6641.625
$ws.Range("L132").Value = 20587.2
$ws.Range("M132").Value = -4111.625
$ws.Range("N132").Value = -25647.2

$ws.Range("H136").Value = 40591.848
$ws.Range("I136").Value = 2414.5715
$ws.Range("J136").Value = 85132
$ws.Range("K136").Value = 7243.7145
$ws.Range("L136").Value = 255396
$ws.Range("M136").Value = -4693.7145
$ws.Range("N136").Value = -260496

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1183.1086
$ws.Range("I5").Value = 950.82855
$ws.Range("J5").Value = 1922.1818
$ws.Range("K5").Value = 2852.48565
$ws.Range("L5").Value = 5766.5454
$ws.Range("M5").Value = -2740.48565
$ws.Range("N5").Value = -5990.5454

$ws.Range("H98").Value = 1582.5
$ws.Range("I98").Value = 899
$ws.Range("J98").Value = 5000
$ws.Range("K98").Value = 2697
$ws.Range("L98").Value = 15000
$ws.Range("M98").Value = -1199
$ws.Range("N98").Value = -17996

$ws.Range("H122").Value = 777.1905
$ws.Range("I122").Value = 292.57144
$ws.Range("J122").Value = 1019.5
$ws.Range("K122").Value = 2633.14296
$ws.Range("L122").Value = 9175.5
$ws.Range("M122").Value = -183.1429600000001
$ws.Range("N122").Value = -14075.5

$ws.Range("H131").Value = 848.3333
$ws.Range("J131").Value = 962.35614
$ws.Range("L131").Value = 2887.06842
$ws.Range("N131").Value = -12967.06842

$ws.Range("H135").Value = 1183.1086
$ws.Range("I135").Value = 950.82855
$ws.Range("J135").Value = 1922.1818
$ws.Range("K135").Value = 8557.45695
$ws.Range("L135").Value = 17299.6362
$ws.Range("M135").Value = -6022.45695
$ws.Range("N135").Value = -22369.6362

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 28160
$ws.Range("J57").Value = 28160
$ws.Range("L57").Value = 28160
$ws.Range("N57").Value = -29800

$ws.Range("H80").Value = 3604.037
$ws.Range("I80").Value = 3365.2144
$ws.Range("K80").Value = 3365.2144
$ws.Range("M80").Value = -2367.2144

$ws.Range("H83").Value = 3604.037
$ws.Range("I83").Value = 3365.2144
$ws.Range("K83").Value = 16826.072
$ws.Range("M83").Value = -11834.072

$ws.Range("H122").Value = 7375
$ws.Range("I122").Value = 9333.333000000001
$ws.Range("K122").Value = 27999.999
$ws.Range("M122").Value = -25549.999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4640
$ws.Range("I7").Value = 3875
$ws.Range("J7").Value = 7700
$ws.Range("K7").Value = 3875
$ws.Range("L7").Value = 7700
$ws.Range("M7").Value = -3763
$ws.Range("N7").Value = -7924

$ws.Range("H16").Value = 683.55554
$ws.Range("J16").Value = 600.1429000000001
$ws.Range("L16").Value = 600.1429000000001
$ws.Range("N16").Value = -940.1429000000001

$ws.Range("H40").Value = 3456.5483
$ws.Range("I40").Value = 3275.8845
$ws.Range("J40").Value = 4396
$ws.Range("K40").Value = 3275.8845
$ws.Range("L40").Value = 4396
$ws.Range("M40").Value = -3139.8845
$ws.Range("N40").Value = -4668

$ws.Range("H100").Value = 1825.7059
$ws.Range("I100").Value = 1204.4445
$ws.Range("K100").Value = 1204.4445
$ws.Range("M100").Value = -663.4445000000001

$ws.Range("H126").Value = 4640
$ws.Range("I126").Value = 3875
$ws.Range("J126").Value = 7700
$ws.Range("K126").Value = 11625
$ws.Range("L126").Value = 23100
$ws.Range("M126").Value = -9155
$ws.Range("N126").Value = -28040

$ws.Range("H132").Value = 2719
$ws.Range("I132").Value = 1986.8334
$ws.Range("J132").Value = 5856.857
$ws.Range("K132").Value = 5960.5002
$ws.Range("L132").Value = 17570.571
$ws.Range("M132").Value = -3430.5002
$ws.Range("N132").Value = -22630.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1643.75
$ws.Range("I81").Value = 267.625
$ws.Range("J81").Value = 4396
$ws.Range("K81").Value = 535.25
$ws.Range("L81").Value = 8792
$ws.Range("M81").Value = 525.75
$ws.Range("N81").Value = -10914

$ws.Range("H84").Value = 1643.75
$ws.Range("I84").Value = 267.625
$ws.Range("J84").Value = 4396
$ws.Range("K84").Value = 2676.25
$ws.Range("L84").Value = 43960
$ws.Range("M84").Value = 2627.75
$ws.Range("N84").Value = -54568

$ws.Range("H100").Value = 586.8570999999999
$ws.Range("I100").Value = 900
$ws.Range("J100").Value = 461.6
$ws.Range("K100").Value = 1800
$ws.Range("L100").Value = 923.2
$ws.Range("M100").Value = -1259
$ws.Range("N100").Value = -2005.2

$ws.Range("H132").Value = 1298.2142
$ws.Range("I132").Value = 516
$ws.Range("K132").Value = 1548
$ws.Range("M132").Value = 982

$ws.Range("H136").Value = 52636860
$ws.Range("I136").Value = 83335860
$ws.Range("K136").Value = 250007580
$ws.Range("M136").Value = -250005030
